$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.210770130157471
$ws.Range("B1").Value = 2.629213094711304
$ws.Range("D1").Value = 2.162957429885864
$ws.Range("E1").Value = 1.16170072555542
